$d = $word.ActiveDocument

$d.Content.Find.Execute("40-27=", $true, $false, $false, $false, $false, $true, 1, $false, "50+0=", 2) | Out-Null
$d.Content.Find.Execute("36+37=", $true, $false, $false, $false, $false, $true, 1, $false, "71-17=", 2) | Out-Null
$d.Content.Find.Execute("94-27=", $true, $false, $false, $false, $false, $true, 1, $false, "66+6=", 2) | Out-Null
$d.Content.Find.Execute("25-18=", $true, $false, $false, $false, $false, $true, 1, $false, "89-71=", 2) | Out-Null
$d.Content.Find.Execute("44-23=", $true, $false, $false, $false, $false, $true, 1, $false, "15+70=", 2) | Out-Null
$d.Content.Find.Execute("8+9=", $true, $false, $false, $false, $false, $true, 1, $false, "4+58=", 2) | Out-Null
$d.Content.Find.Execute("32+48=", $true, $false, $false, $false, $false, $true, 1, $false, "30-27=", 2) | Out-Null
$d.Content.Find.Execute("12+20=", $true, $false, $false, $false, $false, $true, 1, $false, "32-27=", 2) | Out-Null
$d.Content.Find.Execute("56+34=", $true, $false, $false, $false, $false, $true, 1, $false, "64+13=", 2) | Out-Null
$d.Content.Find.Execute("79-43=", $true, $false, $false, $false, $false, $true, 1, $false, "68+23=", 2) | Out-Null
$d.Content.Find.Execute("61-60=", $true, $false, $false, $false, $false, $true, 1, $false, "29+69=", 2) | Out-Null
$d.Content.Find.Execute("85-9=", $true, $false, $false, $false, $false, $true, 1, $false, "6+77=", 2) | Out-Null
$d.Content.Find.Execute("15+26=", $true, $false, $false, $false, $false, $true, 1, $false, "32-8=", 2) | Out-Null
$d.Content.Find.Execute("81-21=", $true, $false, $false, $false, $false, $true, 1, $false, "58+3=", 2) | Out-Null
$d.Content.Find.Execute("72-30=", $true, $false, $false, $false, $false, $true, 1, $false, "69-54=", 2) | Out-Null
$d.Content.Find.Execute("55-32=", $true, $false, $false, $false, $false, $true, 1, $false, "88-31=", 2) | Out-Null
$d.Content.Find.Execute("80+13=", $true, $false, $false, $false, $false, $true, 1, $false, "84-75=", 2) | Out-Null
$d.Content.Find.Execute("91-66=", $true, $false, $false, $false, $false, $true, 1, $false, "82-54=", 2) | Out-Null
$d.Content.Find.Execute("8+40=", $true, $false, $false, $false, $false, $true, 1, $false, "38+19=", 2) | Out-Null
$d.Content.Find.Execute("48+13=", $true, $false, $false, $false, $false, $true, 1, $false, "50+18=", 2) | Out-Null
$d.Content.Find.Execute("48+10=", $true, $false, $false, $false, $false, $true, 1, $false, "0+36=", 2) | Out-Null
$d.Content.Find.Execute("20+30=", $true, $false, $false, $false, $false, $true, 1, $false, "11+76=", 2) | Out-Null
$d.Content.Find.Execute("82-23=", $true, $false, $false, $false, $false, $true, 1, $false, "71-48=", 2) | Out-Null
$d.Content.Find.Execute("89-89=", $true, $false, $false, $false, $false, $true, 1, $false, "40+47=", 2) | Out-Null
$d.Content.Find.Execute("45-1=", $true, $false, $false, $false, $false, $true, 1, $false, "51-1=", 2) | Out-Null
$d.Content.Find.Execute("22+12=", $true, $false, $false, $false, $false, $true, 1, $false, "83-52=", 2) | Out-Null
$d.Content.Find.Execute("2+62=", $true, $false, $false, $false, $false, $true, 1, $false, "19+10=", 2) | Out-Null
$d.Content.Find.Execute("15+36=", $true, $false, $false, $false, $false, $true, 1, $false, "73-0=", 2) | Out-Null
$d.Content.Find.Execute("59+26=", $true, $false, $false, $false, $false, $true, 1, $false, "78+4=", 2) | Out-Null
$d.Content.Find.Execute("26+65=", $true, $false, $false, $false, $false, $true, 1, $false, "29+14=", 2) | Out-Null
$d.Content.Find.Execute("56-11=", $true, $false, $false, $false, $false, $true, 1, $false, "57-32=", 2) | Out-Null
$d.Content.Find.Execute("26-7=", $true, $false, $false, $false, $false, $true, 1, $false, "15+68=", 2) | Out-Null
$d.Content.Find.Execute("12+84=", $true, $false, $false, $false, $false, $true, 1, $false, "43+40=", 2) | Out-Null
$d.Content.Find.Execute("28+52=", $true, $false, $false, $false, $false, $true, 1, $false, "64+30=", 2) | Out-Null
$d.Content.Find.Execute("7+68=", $true, $false, $false, $false, $false, $true, 1, $false, "59+38=", 2) | Out-Null
$d.Content.Find.Execute("47+29=", $true, $false, $false, $false, $false, $true, 1, $false, "28-15=", 2) | Out-Null
$d.Content.Find.Execute("69-28=", $true, $false, $false, $false, $false, $true, 1, $false, "92-91=", 2) | Out-Null
$d.Content.Find.Execute("46+33=", $true, $false, $false, $false, $false, $true, 1, $false, "13+67=", 2) | Out-Null
$d.Content.Find.Execute("78-59=", $true, $false, $false, $false, $false, $true, 1, $false, "97-15=", 2) | Out-Null
$d.Content.Find.Execute("41-6=", $true, $false, $false, $false, $false, $true, 1, $false, "76+5=", 2) | Out-Null
$d.Content.Find.Execute("95-66=", $true, $false, $false, $false, $false, $true, 1, $false, "26-2=", 2) | Out-Null
$d.Content.Find.Execute("7+16=", $true, $false, $false, $false, $false, $true, 1, $false, "35+60=", 2) | Out-Null
$d.Content.Find.Execute("98-34=", $true, $false, $false, $false, $false, $true, 1, $false, "46+22=", 2) | Out-Null
$d.Content.Find.Execute("84-12=", $true, $false, $false, $false, $false, $true, 1, $false, "54-24=", 2) | Out-Null
$d.Content.Find.Execute("9+1=", $true, $false, $false, $false, $false, $true, 1, $false, "9+71=", 2) | Out-Null
$d.Content.Find.Execute("9+35=", $true, $false, $false, $false, $false, $true, 1, $false, "18+8=", 2) | Out-Null
$d.Content.Find.Execute("30+0=", $true, $false, $false, $false, $false, $true, 1, $false, "77-55=", 2) | Out-Null
$d.Content.Find.Execute("78-22=", $true, $false, $false, $false, $false, $true, 1, $false, "72-4=", 2) | Out-Null
$d.Content.Find.Execute("8+70=", $true, $false, $false, $false, $false, $true, 1, $false, "2+41=", 2) | Out-Null
$d.Content.Find.Execute("13+66=", $true, $false, $false, $false, $false, $true, 1, $false, "64+35=", 2) | Out-Null
$d.Content.Find.Execute("8+64=", $true, $false, $false, $false, $false, $true, 1, $false, "95-48=", 2) | Out-Null
$d.Content.Find.Execute("51-17=", $true, $false, $false, $false, $false, $true, 1, $false, "32+18=", 2) | Out-Null
$d.Content.Find.Execute("44+22=", $true, $false, $false, $false, $false, $true, 1, $false, "30+65=", 2) | Out-Null
$d.Content.Find.Execute("44-26=", $true, $false, $false, $false, $false, $true, 1, $false, "58+32=", 2) | Out-Null
$d.Content.Find.Execute("60+12=", $true, $false, $false, $false, $false, $true, 1, $false, "80+10=", 2) | Out-Null
$d.Content.Find.Execute("56+11=", $true, $false, $false, $false, $false, $true, 1, $false, "94-71=", 2) | Out-Null
$d.Content.Find.Execute("84-83=", $true, $false, $false, $false, $false, $true, 1, $false, "86-65=", 2) | Out-Null
$d.Content.Find.Execute("35+28=", $true, $false, $false, $false, $false, $true, 1, $false, "46-14=", 2) | Out-Null
$d.Content.Find.Execute("19+79=", $true, $false, $false, $false, $false, $true, 1, $false, "40-24=", 2) | Out-Null
$d.Content.Find.Execute("12+18=", $true, $false, $false, $false, $false, $true, 1, $false, "53-2=", 2) | Out-Null
$d.Content.Find.Execute("74-42=", $true, $false, $false, $false, $false, $true, 1, $false, "33-11=", 2) | Out-Null
$d.Content.Find.Execute("34+36=", $true, $false, $false, $false, $false, $true, 1, $false, "80-77=", 2) | Out-Null
$d.Content.Find.Execute("86-20=", $true, $false, $false, $false, $false, $true, 1, $false, "7+28=", 2) | Out-Null
$d.Content.Find.Execute("8+89=", $true, $false, $false, $false, $false, $true, 1, $false, "94-7=", 2) | Out-Null
$d.Content.Find.Execute("22+18=", $true, $false, $false, $false, $false, $true, 1, $false, "90-1=", 2) | Out-Null
$d.Content.Find.Execute("11-9=", $true, $false, $false, $false, $false, $true, 1, $false, "66+0=", 2) | Out-Null
$d.Content.Find.Execute("35-33=", $true, $false, $false, $false, $false, $true, 1, $false, "82-4=", 2) | Out-Null
$d.Content.Find.Execute("31+7=", $true, $false, $false, $false, $false, $true, 1, $false, "25+19=", 2) | Out-Null
$d.Content.Find.Execute("84-52=", $true, $false, $false, $false, $false, $true, 1, $false, "22+33=", 2) | Out-Null
$d.Content.Find.Execute("27-3=", $true, $false, $false, $false, $false, $true, 1, $false, "5+79=", 2) | Out-Null
$d.Content.Find.Execute("95-4=", $true, $false, $false, $false, $false, $true, 1, $false, "32-19=", 2) | Out-Null
$d.Content.Find.Execute("13+63=", $true, $false, $false, $false, $false, $true, 1, $false, "15+1=", 2) | Out-Null
$d.Content.Find.Execute("7+6=", $true, $false, $false, $false, $false, $true, 1, $false, "94-24=", 2) | Out-Null
$d.Content.Find.Execute("22-11=", $true, $false, $false, $false, $false, $true, 1, $false, "7+37=", 2) | Out-Null
$d.Content.Find.Execute("59+29=", $true, $false, $false, $false, $false, $true, 1, $false, "28+67=", 2) | Out-Null
$d.Content.Find.Execute("78-37=", $true, $false, $false, $false, $false, $true, 1, $false, "59+38=", 2) | Out-Null
$d.Content.Find.Execute("24+35=", $true, $false, $false, $false, $false, $true, 1, $false, "25+5=", 2) | Out-Null
$d.Content.Find.Execute("69-7=", $true, $false, $false, $false, $false, $true, 1, $false, "48-4=", 2) | Out-Null
$d.Content.Find.Execute("48+4=", $true, $false, $false, $false, $false, $true, 1, $false, "49-45=", 2) | Out-Null
$d.Content.Find.Execute("47+18=", $true, $false, $false, $false, $false, $true, 1, $false, "24+29=", 2) | Out-Null
$d.Content.Find.Execute("51+16=", $true, $false, $false, $false, $false, $true, 1, $false, "94-32=", 2) | Out-Null
$d.Content.Find.Execute("27+40=", $true, $false, $false, $false, $false, $true, 1, $false, "92-22=", 2) | Out-Null
$d.Content.Find.Execute("64+32=", $true, $false, $false, $false, $false, $true, 1, $false, "91-60=", 2) | Out-Null
$d.Content.Find.Execute("92-33=", $true, $false, $false, $false, $false, $true, 1, $false, "3+79=", 2) | Out-Null
$d.Content.Find.Execute("13+86=", $true, $false, $false, $false, $false, $true, 1, $false, "54-0=", 2) | Out-Null
$d.Content.Find.Execute("62-49=", $true, $false, $false, $false, $false, $true, 1, $false, "24+60=", 2) | Out-Null
$d.Content.Find.Execute("38+17=", $true, $false, $false, $false, $false, $true, 1, $false, "92-51=", 2) | Out-Null
$d.Content.Find.Execute("70-19=", $true, $false, $false, $false, $false, $true, 1, $false, "56+21=", 2) | Out-Null
$d.Content.Find.Execute("83-6=", $true, $false, $false, $false, $false, $true, 1, $false, "91-79=", 2) | Out-Null
$d.Content.Find.Execute("93+1=", $true, $false, $false, $false, $false, $true, 1, $false, "97-69=", 2) | Out-Null
$d.Content.Find.Execute("50-19=", $true, $false, $false, $false, $false, $true, 1, $false, "78-69=", 2) | Out-Null
$d.Content.Find.Execute("21+0=", $true, $false, $false, $false, $false, $true, 1, $false, "42+34=", 2) | Out-Null
$d.Content.Find.Execute("63-7=", $true, $false, $false, $false, $false, $true, 1, $false, "6+90=", 2) | Out-Null
$d.Content.Find.Execute("19+68=", $true, $false, $false, $false, $false, $true, 1, $false, "69-24=", 2) | Out-Null
$d.Content.Find.Execute("49-22=", $true, $false, $false, $false, $false, $true, 1, $false, "15+64=", 2) | Out-Null
$d.Content.Find.Execute("39+38=", $true, $false, $false, $false, $false, $true, 1, $false, "13-11=", 2) | Out-Null
$d.Content.Find.Execute("63-30=", $true, $false, $false, $false, $false, $true, 1, $false, "29+39=", 2) | Out-Null
$d.Content.Find.Execute("87-58=", $true, $false, $false, $false, $false, $true, 1, $false, "59-41=", 2) | Out-Null
$d.Content.Find.Execute("48+29=", $true, $false, $false, $false, $false, $true, 1, $false, "27+72=", 2) | Out-Null
$d.Content.Find.Execute("24-16=", $true, $false, $false, $false, $false, $true, 1, $false, "95-75=", 2) | Out-Null
